{"js": "// The document is a single 20-row x 5-column table of arithmetic problems\n// (\"91-86=5\", \"14+8=22\", ...). The edit swaps every one of the 100 answer\n// strings for a new problem, cell-for-cell, in row-major (reading) order.\n//\n// We replace by table position (row, column) rather than by searching for\n// the old text: a few of the new strings are identical to *other* cells'\n// original text (e.g. \"39+53=92\" is both cell #10's original text and\n// cell #41's replacement text), so a sequential find-and-replace on text\n// would become ambiguous partway through. Addressing cells structurally\n// sidesteps that entirely.\nconst newValues = [\n  \"59+28=87\", \"60-18=42\", \"87+5=92\", \"4+87=91\", \"88+7=95\",\n  \"25+69=94\", \"23+69=92\", \"73-65=8\", \"29+43=72\", \"39+53=92\",\n  \"53-47=6\", \"49+33=82\", \"74-9=65\", \"66-8=58\", \"47+37=84\",\n  \"42-9=33\", \"6+37=43\", \"5+48=53\", \"31-9=22\", \"53-49=4\",\n  \"92-35=57\", \"18+19=37\", \"85-58=27\", \"18+45=63\", \"26+15=41\",\n  \"82-3=79\", \"60-26=34\", \"66+28=94\", \"70-65=5\", \"25+46=71\",\n  \"65+16=81\", \"83-35=48\", \"72-57=15\", \"27+35=62\", \"33-25=8\",\n  \"17+35=52\", \"54+8=62\", \"94-36=58\", \"56-17=39\", \"89+7=96\",\n  \"35+59=94\", \"24+67=91\", \"59+36=95\", \"36+6=42\", \"90-65=25\",\n  \"24+7=31\", \"56+36=92\", \"18+58=76\", \"38+6=44\", \"26+8=34\",\n  \"73-26=47\", \"95-66=29\", \"85-58=27\", \"6+45=51\", \"60-59=1\",\n  \"46-8=38\", \"52+9=61\", \"38+9=47\", \"9+72=81\", \"90-7=83\",\n  \"70-11=59\", \"94-46=48\", \"67+6=73\", \"15+38=53\", \"30-13=17\",\n  \"45+48=93\", \"45-9=36\", \"19+65=84\", \"85-68=17\", \"91-72=19\",\n  \"48+18=66\", \"61-54=7\", \"50-37=13\", \"43+49=92\", \"18+33=51\",\n  \"90-19=71\", \"81-75=6\", \"46+16=62\", \"9+82=91\", \"95-26=69\",\n  \"23+48=71\", \"63-26=37\", \"57-49=8\", \"71-59=12\", \"70-33=37\",\n  \"19+36=55\", \"9+34=43\", \"90-71=19\", \"34-8=26\", \"37+14=51\",\n  \"82-8=74\", \"34+27=61\", \"18+56=74\", \"39+5=44\", \"66+18=84\",\n  \"18+34=52\", \"68-49=19\", \"6+27=33\", \"2+79=81\", \"63+8=71\",\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length !== 1) {\n  throw new Error(\"Expected exactly 1 table, found \" + tables.items.length);\n}\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst columnCount = 5;\nlet idx = 0;\nfor (let r = 0; r < table.rowCount; r++) {\n  for (let c = 0; c < columnCount; c++) {\n    const cell = table.getCell(r, c);\n    // Replace the entire cell-body range so the existing run formatting\n    // (font, size) and paragraph properties (left alignment) are kept,\n    // and only the text content changes - mirroring the source diff.\n    const range = cell.body.getRange();\n    range.insertText(newValues[idx], Word.InsertLocation.replace);\n    idx++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document is a single 20-row x 5-column table of arithmetic problems\n# (\"91-86=5\", \"14+8=22\", ...). The edit swaps every one of the 100 answer\n# strings for a new problem, cell-for-cell, in row-major (reading) order.\n#\n# We replace by table position (row, column) rather than Find/Replace on\n# text: a few of the new strings duplicate *other* cells' original text\n# (e.g. \"39+53=92\" is both cell #10's original text and cell #41's\n# replacement text), so a sequential find-and-replace would become\n# ambiguous partway through. Addressing cells structurally sidesteps that.\n$newValues = @(\n  \"59+28=87\",\"60-18=42\",\"87+5=92\",\"4+87=91\",\"88+7=95\",\n  \"25+69=94\",\"23+69=92\",\"73-65=8\",\"29+43=72\",\"39+53=92\",\n  \"53-47=6\",\"49+33=82\",\"74-9=65\",\"66-8=58\",\"47+37=84\",\n  \"42-9=33\",\"6+37=43\",\"5+48=53\",\"31-9=22\",\"53-49=4\",\n  \"92-35=57\",\"18+19=37\",\"85-58=27\",\"18+45=63\",\"26+15=41\",\n  \"82-3=79\",\"60-26=34\",\"66+28=94\",\"70-65=5\",\"25+46=71\",\n  \"65+16=81\",\"83-35=48\",\"72-57=15\",\"27+35=62\",\"33-25=8\",\n  \"17+35=52\",\"54+8=62\",\"94-36=58\",\"56-17=39\",\"89+7=96\",\n  \"35+59=94\",\"24+67=91\",\"59+36=95\",\"36+6=42\",\"90-65=25\",\n  \"24+7=31\",\"56+36=92\",\"18+58=76\",\"38+6=44\",\"26+8=34\",\n  \"73-26=47\",\"95-66=29\",\"85-58=27\",\"6+45=51\",\"60-59=1\",\n  \"46-8=38\",\"52+9=61\",\"38+9=47\",\"9+72=81\",\"90-7=83\",\n  \"70-11=59\",\"94-46=48\",\"67+6=73\",\"15+38=53\",\"30-13=17\",\n  \"45+48=93\",\"45-9=36\",\"19+65=84\",\"85-68=17\",\"91-72=19\",\n  \"48+18=66\",\"61-54=7\",\"50-37=13\",\"43+49=92\",\"18+33=51\",\n  \"90-19=71\",\"81-75=6\",\"46+16=62\",\"9+82=91\",\"95-26=69\",\n  \"23+48=71\",\"63-26=37\",\"57-49=8\",\"71-59=12\",\"70-33=37\",\n  \"19+36=55\",\"9+34=43\",\"90-71=19\",\"34-8=26\",\"37+14=51\",\n  \"82-8=74\",\"34+27=61\",\"18+56=74\",\"39+5=44\",\"66+18=84\",\n  \"18+34=52\",\"68-49=19\",\"6+27=33\",\"2+79=81\",\"63+8=71\"\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables(1)\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\n$idx = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n  for ($c = 1; $c -le $colCount; $c++) {\n    $cell = $t.Cell($r, $c)\n    # Setting Range.Text on the cell's range replaces just the visible\n    # text and keeps the existing run/paragraph formatting (font, size,\n    # alignment) - mirroring the source diff, which only touches <w:t>.\n    $cell.Range.Text = $newValues[$idx]\n    $idx++\n  }\n}\n"}
